# 2021-05-13 샘플다운로드 문서 commit
#
# The original sheet had long, multi-line header labels (with embedded
# example/placeholder text) in row 1, plus a sample data row (row 2).
# This edit shortens the row-1 header labels and moves the placeholder /
# example text that used to be baked into the header strings into proper
# cell comments instead; the now-unnecessary sample data row is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the header labels in row 1 (the long "label\nexample" text is
# replaced by just the label - the example text becomes a comment below).
$ws.Range("A1").Value = "구분"
$ws.Range("C1").Value = "주민등록번호"
$ws.Range("D1").Value = "교육이수번호"
$ws.Range("E1").Value = "경력시작일"
$ws.Range("F1").Value = "경력종료일"
# B1 ("이름") is unchanged.

# Remove the sample data row (row 2: 1 / 홍길동 / 880131-0000000 / ... ).
$ws.Rows(2).Delete() | Out-Null

# Row 1 no longer needs the tall 52.2pt height that fit the old 3-line
# header text - let Excel re-fit it to the (now single-line) content.
$ws.Rows(1).EntireRow.AutoFit() | Out-Null

# Re-attach the descriptive/example text as cell comments instead of part
# of the header strings.
$ws.Range("A1").AddComment("신규 = 1`n경력 = 2") | Out-Null
$ws.Range("C1").AddComment("000000-0000000") | Out-Null
$ws.Range("D1").AddComment("0000000000") | Out-Null
$ws.Range("E1").AddComment("YYYY-MM-DD") | Out-Null
$ws.Range("F1").AddComment("YYYY-MM-DD") | Out-Null

# Leave the selection where the author ended up after editing.
$ws.Range("E8").Select() | Out-Null
